# Specifications_FitMe.xlsx - "Exceptions" sheet update
# - Reorder columns to ID / Class / Text (Class and Text swap places)
# - Add a new exception row for a database-not-responding failure
# - Resize the ID / Class columns to fit their new content

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

# --- Existing data rows: swap the Text and Class columns ----------------
for ($row = 2; $row -le 4; $row++) {
    $classValue = $ws.Cells.Item($row, 2).Value()
    $textValue  = $ws.Cells.Item($row, 3).Value()
    $ws.Cells.Item($row, 3).Value = $classValue
    $ws.Cells.Item($row, 2).Value = $textValue
}

# --- New row: DataBase.cs failure ---------------------------------------
$ws.Cells.Item(5, 1).Value = "0x0003"
$ws.Cells.Item(5, 3).Value = "Database is not responding"
$ws.Cells.Item(5, 2).Value = "DataBase.cs"

# --- Header row: id/text/class -> ID/Class/Text -------------------------
$ws.Cells.Item(1, 2).Value = "Class"
$ws.Cells.Item(1, 3).Value = "Text"
$ws.Cells.Item(1, 1).Value = "ID"

# --- Fit the ID / Class columns to their widened content -----------------
$ws.Columns.Item(1).ColumnWidth = 5.8
$ws.Columns.Item(2).ColumnWidth = 12.3

# --- Reset selection back to the top of the sheet ------------------------
$ws.Range("A1").Select() | Out-Null
